# Sprint Review Protocol 3 - fill in sprint task rows 6-8 (B21:D23)
# and move the active selection to C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (task 6): "Shooting mergen", Estimate=2, Real=3
$ws.Range("B21").Value = "Shooting mergen"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 3

# Row 22 (task 7): "Bounce mergen", Estimate=2, Real=2
$ws.Range("B22").Value = "Bounce mergen"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 2

# Row 23 (task 8): "Hit registration mergen", Estimate=2, Real=1
$ws.Range("B23").Value = "Hit registration mergen"
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1

# Move the selection the author ended up on
$ws.Range("C23").Select()
